$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates. Some new values (e.g. "1.000", "13.40") look like
# plain decimal numbers, so Excel would silently coerce them to numeric cells
# and drop the significant trailing digits. Force Text format for the write,
# then restore the default "Normal" style so the cell format matches the rest
# of the sheet (only the text content changes, like in the source data).
$ws.Range("D2").Value = "23.900.01"
$ws.Range("D3").Value = "1.648.78"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3893"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3840"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.345"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08439"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.018"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.905"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Value = "1.652.65"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06968"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "23.884.49"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.441"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.906"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.382"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "137.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.722"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.485"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "1.829.03"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08114"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9911"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02924"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.696"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2686"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09122"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7547"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.421"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6936"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.440"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.099"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.224"
$ws.Range("D51").Style = "Normal"

# Column E (Volume 1h) updates
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("E20").Value = "  -2.85%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("E25").Value = "  -2.98%  "
$ws.Range("E26").Value = "  -5.12%  "
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("E30").Value = "  -2.12%  "
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("E32").Value = "  -1.58%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("E35").Value = "  -4.57%  "
$ws.Range("E36").Value = "  -3.45%  "
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("E39").Value = "  -4.44%  "
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("E51").Value = "  -0.73%  "
